$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Shift old rows 453:548 down by 2 rows, into new rows 455:550.
# Use Copy with a destination so that both values and formatting (e.g. date
# number format on column D) are carried along in one operation.
$srcBlock = $ws.Range("A453:R548")
$srcBlock.Copy($ws.Range("A455"))

# Step 2: New row 454 takes the values that row 452 had *before* this edit
# (i.e. the data point that used to live at row 452).
$ws.Cells.Item(454, 1).Value2  = 9
$ws.Cells.Item(454, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(454, 3).Value2  = "Metropolitana"
$ws.Cells.Item(454, 4).Value2  = 45005
$ws.Cells.Item(454, 5).Value2  = 13
$ws.Cells.Item(454, 6).Value2  = 100112052
$ws.Cells.Item(454, 7).Value2  = "Albahaca"
$ws.Cells.Item(454, 8).Value2  = "Sin especificar"
$ws.Cells.Item(454, 9).Value2  = "Primera"
$ws.Cells.Item(454, 10).Value2 = 340
$ws.Cells.Item(454, 11).Value2 = 3000
$ws.Cells.Item(454, 12).Value2 = 3500
$ws.Cells.Item(454, 13).Value2 = 3250
$ws.Cells.Item(454, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(454, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(454, 16).Value2 = 542
$ws.Cells.Item(454, 17).Value2 = 6
$ws.Cells.Item(454, 18).Value2 = "Hortaliza"

# Step 3: New row 453 is a brand-new data point.
$ws.Cells.Item(453, 1).Value2  = 9
$ws.Cells.Item(453, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(453, 3).Value2  = "Metropolitana"
$ws.Cells.Item(453, 4).Value2  = 45015
$ws.Cells.Item(453, 5).Value2  = 13
$ws.Cells.Item(453, 6).Value2  = 100112052
$ws.Cells.Item(453, 7).Value2  = "Albahaca"
$ws.Cells.Item(453, 8).Value2  = "Sin especificar"
$ws.Cells.Item(453, 9).Value2  = "Segunda"
$ws.Cells.Item(453, 10).Value2 = 160
$ws.Cells.Item(453, 11).Value2 = 2800
$ws.Cells.Item(453, 12).Value2 = 2800
$ws.Cells.Item(453, 13).Value2 = 2800
$ws.Cells.Item(453, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(453, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(453, 16).Value2 = 467
$ws.Cells.Item(453, 17).Value2 = 6
$ws.Cells.Item(453, 18).Value2 = "Hortaliza"

# Step 4: Row 452 keeps most of its data, but the date and a few price
# fields are corrected.
$ws.Cells.Item(452, 4).Value2  = 45015
$ws.Cells.Item(452, 12).Value2 = 3000
$ws.Cells.Item(452, 13).Value2 = 3000
$ws.Cells.Item(452, 16).Value2 = 500
